# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.140.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.50%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.656.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.47%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.53%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'218.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.11%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.5244"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.54%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.48%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2620"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.06293"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'20.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.27%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07801"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.31%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.502"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.34%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.679.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.68%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.883.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.41%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.5549"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.37%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -2.34%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -0.89%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'26.152.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.49%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.005"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.50%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'4.637"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.88%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'195.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.18%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E23").Value = "'  -1.13%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -0.49%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'146.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "'  -1.85%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'7.167"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.25%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.76%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.79%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.05727"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.44%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.270"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.63%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.489"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.65%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.347"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.15%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.588"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.30%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.804"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.64%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.9524"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.83%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.23%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.5699"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.69%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.36%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'5.957"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.11%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.059.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.23%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.8443"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.21%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.005"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.47%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'103.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.58%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.794.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.41%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'57.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.62%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'BabyDogeCoin"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.0₈107"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.10%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'Cronos"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'0.05398"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +4.56%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.009"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.38%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.4401"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.46%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'8.010"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.52%  "
$ws.Range("E51").Style = "Normal"
